$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.581.25"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "3.270.95"
$ws.Range("E3").Value = "  -5.40%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'592.20"
$ws.Range("E5").Value = "  -3.10%  "

$ws.Range("D6").Value = "'150.91"
$ws.Range("E6").Value = "  -9.69%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.262.76"
$ws.Range("E8").Value = "  -5.60%  "

$ws.Range("D9").Value = "'0.544"
$ws.Range("E9").Value = "  -8.34%  "

$ws.Range("E10").Value = "  -10.55%  "

$ws.Range("E11").Value = "  -4.80%  "

$ws.Range("D12").Value = "'0.506"
$ws.Range("E12").Value = "  -10.27%  "

$ws.Range("D13").Value = "'38.53"
$ws.Range("E13").Value = "  -13.03%  "

$ws.Range("E14").Value = "  -8.21%  "

$ws.Range("D15").Value = "3.786.77"
$ws.Range("E15").Value = "  -5.93%  "

$ws.Range("D16").Value = "67.507.30"
$ws.Range("E16").Value = "  -3.47%  "

$ws.Range("D17").Value = "3.268.90"
$ws.Range("E17").Value = "  -5.63%  "

$ws.Range("E18").Value = "  -5.28%  "

$ws.Range("D19").Value = "'533.06"
$ws.Range("E19").Value = "  -8.78%  "

$ws.Range("D20").Value = "'7.15"
$ws.Range("E20").Value = "  -12.60%  "

$ws.Range("D21").Value = "'15.02"
$ws.Range("E21").Value = "  -12.45%  "

$ws.Range("D22").Value = "'0.759"

$ws.Range("E23").Value = "  -12.27%  "

$ws.Range("D24").Value = "'85.62"
$ws.Range("E24").Value = "  -10.30%  "

$ws.Range("D25").Value = "'13.59"
$ws.Range("E25").Value = "  -10.17%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("E27").Value = "  -10.17%  "

$ws.Range("D28").Value = "'8.13"
$ws.Range("E28").Value = "  -5.75%  "

$ws.Range("D29").Value = "'2.16"
$ws.Range("E29").Value = "  -12.04%  "

$ws.Range("D30").Value = "'29.28"
$ws.Range("E30").Value = "  -10.90%  "

$ws.Range("D31").Value = "'2.70"
$ws.Range("E31").Value = "  -4.47%  "

$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  -6.57%  "

$ws.Range("D33").Value = "'6.64"
$ws.Range("E33").Value = "  -15.43%  "

$ws.Range("D34").Value = "'5.75"
$ws.Range("E34").Value = "  -12.33%  "

$ws.Range("D35").Value = "'518.80"
$ws.Range("E35").Value = "  -11.54%  "

$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").Value = "'0.0448"
$ws.Range("E37").Value = "  -6.45%  "

$ws.Range("D38").Value = "'53.44"
$ws.Range("E38").Value = "  -4.79%  "

$ws.Range("D39").Value = "'0.0859"
$ws.Range("E39").Value = "  -10.45%  "

$ws.Range("D40").Value = "'9.01"
$ws.Range("E40").Value = "  -15.04%  "

$ws.Range("D41").Value = "'0.126"
$ws.Range("E41").Value = "  -10.14%  "

$ws.Range("D42").Value = "'2.80"
$ws.Range("E42").Value = "  -10.74%  "

$ws.Range("D43").Value = "2.953.28"
$ws.Range("E43").Value = "  -8.87%  "

$ws.Range("D44").Value = "'0.267"
$ws.Range("E44").Value = "  -9.46%  "

$ws.Range("D45").Value = "0.0₃0590"
$ws.Range("E45").Value = "  -15.29%  "

$ws.Range("D46").Value = "'2.20"
$ws.Range("E46").Value = "  -8.34%  "

$ws.Range("D47").Value = "'26.84"
$ws.Range("E47").Value = "  -12.67%  "

$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").Value = "'2.33"
$ws.Range("E49").Value = "  -15.67%  "

$ws.Range("E50").Value = "  -9.41%  "

$ws.Range("D51").Value = "'123.75"
$ws.Range("E51").Value = "  -7.47%  "
